# "Final touches in slides"
# Reposition a handful of caption textboxes and a picture, widen the
# "Report" label (and rename it to "Report/webpage").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# TextBox 37 ("Documentation") -> nudge left/up slightly
$docBox = $s.Shapes.Item(12)
$docBox.Left = 317.4667816535433
$docBox.Top  = 391.53639295275593

# TextBox 38 ("Figures") -> nudge left/down slightly
$figBox = $s.Shapes.Item(13)
$figBox.Left = 321.0104834409449
$figBox.Top  = 308.7188976377953

# TextBox 39 ("Code") -> nudge left/down slightly
$codeBox = $s.Shapes.Item(14)
$codeBox.Left = 321.0104834409449
$codeBox.Top  = 224.32251968503937

# TextBox 42 ("Report") -> widen box and update label text
$reportBox = $s.Shapes.Item(17)
$reportBox.TextFrame.TextRange.Text = "Report/webpage"
$reportBox.Width = 141.92252368503935

# Picture 63 -> move to the right, slightly up
$pic63 = $s.Shapes.Item(30)
$pic63.Left = 498.2596062992126
$pic63.Top  = 442.42858267716537
